# Updated cryptos list on Wed Jun  7 05:59:36 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto listing sheet, and re-sorts the two coins (Toncoin /
# EthereumClassic) whose ranking swapped positions between rows 28 and 29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values are plain text (e.g. thousand-dot formatted), so force the
# cell to Text ("@") number format before assigning, keeping Excel from
# re-interpreting the string as a numeric value and dropping formatting
# such as trailing zeros.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.911.55"
$ws.Range("E2").Value = "  +4.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.880.51"
$ws.Range("E3").Value = "  +3.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "278.68"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5326"
$ws.Range("E7").Value = "  +4.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3461"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06996"
$ws.Range("E10").Value = "  +5.01%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8083"
$ws.Range("E12").Value = "  -2.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07786"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.871.73"
$ws.Range("E14").Value = "  +3.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.94"
$ws.Range("E15").Value = "  +3.70%  "
$ws.Range("E16").Value = "  +2.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.59"
$ws.Range("E17").Value = "  +3.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008052"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.981.47"
$ws.Range("E21").Value = "  +4.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.110.48"
$ws.Range("E22").Value = "  +3.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.755"
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.206"
$ws.Range("E25").Value = "  +2.05%  "
$ws.Range("E26").Value = "  +8.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.87"
$ws.Range("E27").Value = "  +3.70%  "

# Rows 28 and 29 swapped ranking order (Toncoin now ranks above
# EthereumClassic), so update the Coin name, Link, Price and
# Volume(1h) columns for both rows accordingly.
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.666"
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.37"
$ws.Range("E29").Value = "  +1.69%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "113.73"
$ws.Range("E30").Value = "  +3.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.372"
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.338"
$ws.Range("E32").Value = "  +2.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08914"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04948"
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("E35").Value = "  +4.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7348"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.886"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.291"
$ws.Range("E38").Value = "  +4.41%  "
$ws.Range("E39").Value = "  +1.32%  "
$ws.Range("E40").Value = "  +0.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5166"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9619"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.17"
$ws.Range("E43").Value = "  +4.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.208"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.134"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4525"
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1350"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.353"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.33"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("E51").Value = "  +2.06%  "
